$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds header "K" (formerly derived from "Strike#").
# Regenerated values per commit message: "use K instead of Strike#, regen std/mean, calc and write s_vals"
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 2
